# This workbook's input sheets were authored for the "Weekly" payroll run
# but were actually meant to drive the "Monthly" tax process, and the
# "DO NOT TOUCH" employee placeholder was pointing at the wrong employee
# number. Fix both: rename the two worksheets and correct every cell that
# references them / the employee marker.

$wb = $excel.ActiveWorkbook

# --- Rename worksheets (Weekly -> Monthly) ---
$wsGeneral = $wb.Worksheets.Item(2)
$wsGeneral.Name = "GeneralTaxRateMonthly"

$wsProcess = $wb.Worksheets.Item(3)
$wsProcess.Name = "ProcessPayrollForMonthlyTax"

$wsFirst = $wb.Worksheets.Item(1)
$wsReports = $wb.Worksheets.Item(4)

# --- Sheet "first": update the script-name references to match the rename ---
$wsFirst.Range("A3").Value = "GeneralTaxRateMonthly"
$wsFirst.Range("A4").Value = "ProcessPayrollForMonthlyTax"

# --- Fix "DO NOT TOUCH AUTOMATION EMP 107" -> "... EMP 105" on every sheet ---
$wsGeneral.Range("A2").Value = "DO NOT TOUCH AUTOMATION EMP 105"
$wsProcess.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 105"
$wsReports.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 105"

# --- Update the saved selection / active cell on each sheet ---
$wsFirst.Activate() | Out-Null
$wsFirst.Range("F5").Select() | Out-Null

$wsGeneral.Activate() | Out-Null
$wsGeneral.Range("D10").Select() | Out-Null

$wsProcess.Activate() | Out-Null
$wsProcess.Range("F11").Select() | Out-Null

$wsReports.Activate() | Out-Null
$wsReports.Range("K8").Select() | Out-Null
